$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Mars/Avril" consumption for 2017 (row 4, column C)
$ws.Range("C4").Value = 107.55

# Move the active selection to E14 (matches the saved selection state in the file)
$ws.Range("E14").Select()
